$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates ---
# "Volume 29   Number  50" -> "Volume 29   Number  51"
$ws.Range("A8").Value = "Volume 29   Number  51"

# "Report Covering the Week  12/12/2022  Through  12/18/2022"
# -> "Report Covering the Week  12/19/2022  Through  12/25/2022"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# --- Row 15 (Rape) ---
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 13
$ws.Range("K15").Value = -27.777777777777
$ws.Range("L15").Value = -23.529411764705

# --- Row 16 (Robbery) ---
$ws.Range("F16").Value = 11
$ws.Range("H16").Value = 120
$ws.Range("I16").Value = 82
$ws.Range("J16").Value = 65
$ws.Range("K16").Value = 26.153846153846
$ws.Range("L16").Value = 3.797468354430

# --- Row 17 (Fel. Assault) ---
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -12.5
$ws.Range("I17").Value = 211
$ws.Range("J17").Value = 172
$ws.Range("K17").Value = 22.674418604651
$ws.Range("L17").Value = 11.052631578947

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = "0"
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 20
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = 11.25
$ws.Range("L18").Value = 3.488372093023

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = -42.857142857142
$ws.Range("I19").Value = 391
$ws.Range("J19").Value = 331
$ws.Range("K19").Value = 18.126888217522
$ws.Range("L19").Value = 29.042904290429

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 10
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 113
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = 109.259259259259
$ws.Range("L20").Value = 109.259259259259

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = -15.384615384615
$ws.Range("I21").Value = 900
$ws.Range("J21").Value = 723
$ws.Range("K21").Value = 24.481327800829
$ws.Range("L21").Value = 22.282608695652

# --- Row 23 (Housing) ---
$ws.Range("F23").Value = 1

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -47.368421052631
$ws.Range("F24").Value = 94
$ws.Range("G24").Value = 119
$ws.Range("H24").Value = -21.008403361344
$ws.Range("I24").Value = 1376
$ws.Range("J24").Value = 965
$ws.Range("K24").Value = 42.590673575129
$ws.Range("L24").Value = 58.160919540229

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 16.666666666666
$ws.Range("I25").Value = 485
$ws.Range("J25").Value = 404
$ws.Range("K25").Value = 20.049504950495
$ws.Range("L25").Value = 38.571428571428

# --- Row 26 (UCR Rape*) ---
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 23
$ws.Range("K26").Value = -23.333333333333
$ws.Range("L26").Value = -4.166666666666

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 42
$ws.Range("K27").Value = -8.695652173913
$ws.Range("L27").Value = 82.608695652173
